# Add the new USER_REFINED_PK field as a new trailing column (X) on the header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell X1 -> "USER_REFINED_PK" (creates/reuses the shared string,
# and grows the sheet dimension + row span automatically).
$ws.Range("X1").Value = "USER_REFINED_PK"

# Match the width used by the sheet's other header columns for the new column X.
$ws.Columns.Item(24).ColumnWidth = 18

# Move the selection to reflect where the user ended up after adding the column.
$ws.Range("Y5").Select()
